# Update "想去人数" (number of interested people) values in the
# "展览" (Exhibitions) and "全部类型" (All types) sheets, reflecting the
# refreshed data output at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 704
$wsExhibit.Range("F5").Value = 2334
$wsExhibit.Range("F6").Value = 49
$wsExhibit.Range("F7").Value = 3508
$wsExhibit.Range("F8").Value = 464
$wsExhibit.Range("F9").Value = 884

# --- Sheet "全部类型" ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 704
$wsAll.Range("F6").Value = 2334
$wsAll.Range("F7").Value = 49
$wsAll.Range("F8").Value = 3508
$wsAll.Range("F9").Value = 464
$wsAll.Range("F10").Value = 885
